# This script applies a cyclic rotation of data across rows 3, 4 and 5
# on the "Artfynd" worksheet:
#   old row 3 values -> new row 4
#   old row 4 values -> new row 5
#   old row 5 values -> new row 3
#
# Only columns A, B, D, E, F, G, H, P, Q, R change; all other columns
# already hold identical values across these three rows, so they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

# Capture the current ("before") values for the three affected rows.
$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("${col}3").Value()
    $row4[$col] = $ws.Range("${col}4").Value()
    $row5[$col] = $ws.Range("${col}5").Value()
}

# Write the rotated values: row5 -> row3, row3 -> row4, row4 -> row5.
foreach ($col in $cols) {
    $ws.Range("${col}3").Value = $row5[$col]
    $ws.Range("${col}4").Value = $row3[$col]
    $ws.Range("${col}5").Value = $row4[$col]
}
